$d = $word.ActiveDocument

# This edit is a textually no-op change: each targeted run is "replaced"
# with identical text via Find/Replace so Word re-serializes the run's
# <w:t> element, dropping the now-unnecessary xml:space="preserve" on runs
# whose text has no leading/trailing whitespace (matching the upstream diff).

function Replace-InParagraph($paraIndex, [string]$text) {
    $rng = $d.Paragraphs($paraIndex).Range
    $rng.Find.Execute($text, $true, $true, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null
}

# For runs that carry no direct character formatting (an empty <w:rPr/>),
# touching their text regenerates the run and silently drops the empty
# <w:rPr/> element. Nudging a character-formatting property (set then
# unset) on the *exact* run range afterwards forces the engine to re-emit
# the (still empty) <w:rPr/> without touching already-fixed text or
# neighboring runs/paragraph marks.
function Restore-EmptyRPr($paraIndex, [string]$text) {
    $p = $d.Paragraphs($paraIndex).Range
    $found = $d.Range($p.Start, $p.End)
    $found.Find.Execute($text, $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $found.Bold = 1
    $found.Bold = 0
}

# Paragraph 1: "In the quiet depths of the soul, a war is" / "where soldiers, memories, and fears gather" / "preparing for internal conflict."
Replace-InParagraph 1 "In the quiet depths of the soul, a war is"
Replace-InParagraph 1 "where soldiers, memories, and fears gather"
Replace-InParagraph 1 "preparing for internal conflict."

# Paragraph 3: "However, the shadows of the past cast long, unsettling shades" / "each a ghostly reminder of wounds yet to heal."
Replace-InParagraph 3 "However, the shadows of the past cast long, unsettling shades"
Replace-InParagraph 3 "each a ghostly reminder of wounds yet to heal."

# Paragraph 5: "leaving the battlefield of the mind in a state of perpetual uncertainty" / "."
Replace-InParagraph 5 "leaving the battlefield of the mind in a state of perpetual uncertainty"

# The lone "." run has no direct formatting and can't be matched with
# MatchWholeWord, so target just the final character of the paragraph
# (the run text is only ".").
$p5 = $d.Paragraphs(5).Range
$lastChar = $d.Range($p5.End - 2, $p5.End - 1)
$lastChar.Find.Execute(".", $true, $false, $false, $false, $false, $true, 1, $false, ".", 2) | Out-Null
$p5b = $d.Paragraphs(5).Range
$lastChar2 = $d.Range($p5b.End - 2, $p5b.End - 1)
$lastChar2.Bold = 1
$lastChar2.Bold = 0

# Paragraph 9: "colors" (no direct formatting)
Replace-InParagraph 9 "colors"
Restore-EmptyRPr 9 "colors"

# Paragraph 13: "The heart, torn..." (no direct formatting)
Replace-InParagraph 13 "The heart, torn between the ideologies of love and regret, finds itself in a dance of emotional conflict, with the tender caresses of love clashing against the relentless grasp of regret, turning the heart into a battleground of emotions."
Restore-EmptyRPr 13 "The heart, torn between the ideologies of love and regret, finds itself in a dance of emotional conflict, with the tender caresses of love clashing against the relentless grasp of regret, turning the heart into a battleground of emotions."

# Paragraph 15: "Amid the chaos..." (no direct formatting)
Replace-InParagraph 15 "Amid the chaos, the diplomacy of reason and emotion negotiates in the war room of consciousness, their debates echoing in the chambers of rationality, striving to find a ceasefire amidst the din of internal discord."
Restore-EmptyRPr 15 "Amid the chaos, the diplomacy of reason and emotion negotiates in the war room of consciousness, their debates echoing in the chambers of rationality, striving to find a ceasefire amidst the din of internal discord."

# Paragraph 19: "In the ensuing silence..." (no direct formatting)
Replace-InParagraph 19 "In the ensuing silence, the mind becomes a transformed battlefield—a landscape scarred by the echoes of internal struggles, testifying to the resilience of the human spirit in navigating the complexities of the psyche."
Restore-EmptyRPr 19 "In the ensuing silence, the mind becomes a transformed battlefield—a landscape scarred by the echoes of internal struggles, testifying to the resilience of the human spirit in navigating the complexities of the psyche."
